$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15, shifting rows 15-53 down to 16-54.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44487
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100112040
$ws.Cells.Item(15, 7).Value = "Cilantro"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 1000
$ws.Cells.Item(15, 12).Value = 1200
$ws.Cells.Item(15, 13).Value = 1100
$ws.Cells.Item(15, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 550
$ws.Cells.Item(15, 17).Value = 2
$ws.Cells.Item(15, 18).Value = "Hortaliza"
